# Excel_De05.xlsx update: remove old "Câu hỏi 3" (thuế xuất) question block from the
# CAUHOI sheet and renumber the remaining questions, bumping the points value for
# what becomes the new "Câu hỏi 4".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CAUHOI")

# Delete the two rows holding the old Question 3 ("Điền giá trị cho cột thuế xuất...")
# and the blank spacer row that followed it. Everything below shifts up by two rows.
$ws.Rows("7:8").Delete()

# Re-write the renumbered question text for the rows that shifted into place.
$ws.Range("A7").Value = "Câu hỏi 3. (0,5đ) Thuế được tính là 10% nếu những mặt hàng được xuất trong vòng 15 ngày kể từ ngày nhập. Còn lại là 5%"
$ws.Range("A9").Value = "Câu hỏi 4. (1,0 đ) Tính tổng thuế của các Mã hàng (M, R, C) và điền vào Bảng 1"
$ws.Range("A11").Value = "Câu hỏi 5. (0,5 đ) Vẽ đồ thị hình cột (Cluster Column) cho tổng thuế các Mã hàng ở Bảng 1"
$ws.Range("A13").Value = "Câu hỏi 6. (0,5đ) Sắp xếp bảng tính tăng dần theo Nhà cung cấp. Nếu có cùng Nhà cung cấp thì sắp xếp giảm dần theo Ngày nhập"
$ws.Range("A15").Value = "Câu hỏi 7 (0,5đ) Lọc ra danh sách các nhà cung cấp là Ba Sao hoặc ngày nhập trước "

# Match the saved selection state from the edited workbook.
$ws.Range("A16").Select()
